# Update "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps as part of regenerating the
# Handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the first file.
$wsOverview.Range("G2").Value = "2016-08-27 01:05:15"

# zh-cn sheet: Correspond Handoff / Handback Datetime for the first file.
$wsZhCn.Range("H2").Value = "2016-08-27 01:05:11"
$wsZhCn.Range("K2").Value = "2016-08-27 01:05:28"

# de-de sheet: Correspond Handoff / Handback Datetime for the first file.
# (H2 shares the same underlying text as Overview!G2.)
$wsDeDe.Range("H2").Value = "2016-08-27 01:05:15"
$wsDeDe.Range("K2").Value = "2016-08-27 01:05:35"
